# Updates the cryptos list (Coin price / 1h-volume columns, and two
# row-pairs whose rank order swapped) to match the latest scrape.
# Generated for commit: "Updated cryptos list ... with GitHub Actions"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=2; D="54.061.55"; E="  -4.29%  "},
    @{Row=3; D="2.282.83"; E="  -3.89%  "},
    @{Row=4; E="  +0.38%  "},
    @{Row=5; D="490.15"; E="  -2.74%  "},
    @{Row=6; D="127.48"; E="  -3.49%  "},
    @{Row=7; E="  +0.45%  "},
    @{Row=8; D="0.525"; E="  -4.47%  "},
    @{Row=9; D="2.284.07"; E="  -3.94%  "},
    @{Row=10; D="0.0935"; E="  -5.26%  "},
    @{Row=11; E="  -1.20%  "},
    @{Row=12; D="4.74"; E="  +1.93%  "},
    @{Row=13; D="0.315"; E="  -4.43%  "},
    @{Row=14; D="2.728.08"; E="  -2.47%  "},
    @{Row=15; D="21.23"; E="  -1.32%  "},
    @{Row=16; D="54.094.08"; E="  -4.11%  "},
    @{Row=17; D="0.0000128"; E="  -3.24%  "},
    @{Row=18; D="2.293.68"; E="  -1.70%  "},
    @{Row=19; D="9.65"; E="  -3.98%  "},
    @{Row=20; E="  -1.18%  "},
    @{Row=21; D="302.68"; E="  -1.99%  "},
    @{Row=22; D="6.16"; E="  -0.88%  "},
    @{Row=24; D="63.78"; E="  -2.10%  "},
    @{Row=25; E="  +0.22%  "},
    @{Row=26; D="0.366"; E="  -1.65%  "},
    @{Row=27; E="  -4.79%  "},
    @{Row=28; D="7.07"; E="  -2.51%  "},
    @{Row=29; D="169.00"; E="  -2.02%  "},
    @{Row=30; D="0.0₃0696"; E="  -3.33%  "},
    @{Row=31; D="1.61"; E="  -2.05%  "},
    @{Row=32; D="1.00"; E="  +0.05%  "},
    @{Row=33; B="FirstDigitalUSD"; C="https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"; D="1.00"; E="  +0.58%  "},
    @{Row=34; B="Aptos"; C="https://coinranking.com/coin/HGYj5JCv5+aptos-apt"; D="5.74"; E="  -0.25%  "},
    @{Row=35; D="1.05"; E="  -4.87%  "},
    @{Row=36; D="17.55"; E="  -1.00%  "},
    @{Row=37; D="1.16"; E="  -2.47%  "},
    @{Row=38; D="0.846"; E="  +5.23%  "},
    @{Row=39; D="3.60"; E="  -5.26%  "},
    @{Row=40; D="35.72"; E="  -1.05%  "},
    @{Row=41; D="0.367"; E="  -1.27%  "},
    @{Row=42; E="  -3.70%  "},
    @{Row=43; D="3.31"; E="  -1.68%  "},
    @{Row=44; D="123.42"; E="  -5.78%  "},
    @{Row=45; D="4.66"; E="  -3.05%  "},
    @{Row=46; D="0.0881"; E="  -2.99%  "},
    @{Row=47; D="0.542"; E="  -3.79%  "},
    @{Row=48; D="236.50"; E="  -3.78%  "},
    @{Row=49; D="0.0472"; E="  -1.99%  "},
    @{Row=50; B="EnergySwap"; C="https://coinranking.com/coin/SbWqqTui-+energyswap-ens"; D="16.44"; E="  -2.76%  "},
    @{Row=51; B="VeChain"; C="https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"; D="0.0202"; E="  -3.01%  "}
)

foreach ($u in $updates) {
    if ($u.ContainsKey('B')) {
        $ws.Cells.Item($u.Row, 2).Value = $u.B
    }
    if ($u.ContainsKey('C')) {
        $ws.Cells.Item($u.Row, 3).Value = $u.C
    }
    if ($u.ContainsKey('D')) {
        # Column D ("Price") holds plain-text, locale-formatted numbers
        # (e.g. thousand-dot groups like "54.061.55"). Force the cell to
        # Text format first so Excel doesn't reinterpret/round numeric-
        # looking values (e.g. "490.15" -> 490.14999999999998) or strip
        # trailing zeros (e.g. "169.00" -> 169).
        $cell = $ws.Cells.Item($u.Row, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
    }
    if ($u.ContainsKey('E')) {
        $ws.Cells.Item($u.Row, 5).Value = $u.E
    }
}

Write-Output "Updated $($updates.Count) rows"
